$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44: E44 set values reordered in the set representation
$ws.Range("E44").Value = "{'list', 'str'}"

# Row 45: E45 changes from 'str' to 'list', and F45 changes from Loss (red) to Neutral (orange)
$ws.Range("E45").Value = "list"
$ws.Range("F45").Value = "Neutral"
$ws.Range("F45").Interior.Color = 42495

# Row 57: D57 (PyType Wins) changes from 5 to 4
$ws.Range("D57").Value = 4

# Row 58: "Scalpel Accuracy:" label moves from C58 to E58, and value moves from D58 (1000) to F58 (92.73)
$ws.Range("C58").ClearContents()
$ws.Range("D58").ClearContents()
$ws.Range("E58").Value = "Scalpel Accuracy:"
$ws.Range("F58").Value = 92.73

# Row 59: label text change and accuracy value change
$ws.Range("E59").Value = "Accuracy vs PyType"
$ws.Range("F59").Value = 25
